$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price values so Excel does not
# auto-coerce them to numbers (they must stay text, matching the source data).
$textCells = "D5","D6","D14","D22","D23","D30","D32","D37","D48","D51"
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.690.68'
$ws.Range('D3').Value = '3.443.67'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '579.81'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').Value = '147.96'
$ws.Range('E6').Value = '  +0.39%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +0.32%  '
$ws.Range('E9').Value = '  +4.08%  '
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('E11').Value = '  +2.57%  '
$ws.Range('E12').Value = '  -1.09%  '
$ws.Range('E13').Value = '  +2.09%  '
$ws.Range('D14').Value = '28.23'
$ws.Range('E14').Value = '  -5.53%  '
$ws.Range('E15').Value = '  -1.42%  '
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').Value = '62.715.58'
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('E20').Value = '  -3.00%  '
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').Value = '0.561'
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '75.28'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D25').Value = '3.581.19'
$ws.Range('E25').Value = '  -1.30%  '
$ws.Range('E26').Value = '  -1.70%  '
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = '7.98'
$ws.Range('E30').Value = '  -3.26%  '
$ws.Range('E31').Value = '  -1.31%  '
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('E33').Value = '  -4.44%  '
$ws.Range('E34').Value = '  -2.58%  '
$ws.Range('E35').Value = '  +3.21%  '
$ws.Range('E36').Value = '  +0.54%  '
$ws.Range('D37').Value = '31.88'
$ws.Range('E37').Value = '  +1.16%  '
$ws.Range('E38').Value = '  -2.00%  '
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('E40').Value = '  -1.05%  '
$ws.Range('E41').Value = '  +0.33%  '
$ws.Range('E42').Value = '  -2.47%  '
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('E45').Value = '  -2.50%  '
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('D47').Value = '2.567.52'
$ws.Range('E47').Value = '  -1.83%  '
$ws.Range('D48').Value = '6.92'
$ws.Range('E48').Value = '  +1.93%  '
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('E50').Value = '  -3.89%  '
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.11%  '
